# Apply weekly report refresh: update generated timestamp and zero out
# billed amounts / pricing figures (no-violation / re-run scenario).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report generation timestamp (D5)
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"

# Total Billed Amount (summary box)
$ws.Range("C8").Value = 0

# Per-line pricing values zeroed out
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0

# TOTAL row
$ws.Range("H25").Value = 0
